$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values with repulled data
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -9
$ws.Range("F4").Value = -7
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -1
